# Important critical updates - need to RAISE INVOICE for 50%
#
# Adds a new "Revision-10" line item (row 28) to the schedule table,
# pushing the totals row (and the two blank spacer rows below it) down
# by one row, and updates the running total formula/value accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 28. This shifts:
#      old row 28 (bold totals row)      -> row 29
#      old row 31 (blank spacer row)     -> row 32
#      old row 39 (blank spacer row)     -> row 40
#    while preserving each shifted row's existing per-cell formatting.
$ws.Rows.Item(28).Insert()

# 2. The freshly inserted row 28 has no formatting yet - clone the
#    look of the row directly above it (row 27, the last real data row)
#    so the new row matches the rest of the table.
$ws.Range("B27:J27").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Fill in the new line item's data.
$ws.Range("B28").Value = 22
$ws.Range("C28").Value = "Revision-10_File duplicate checks & Date related"
$ws.Range("D28").Value = "13 - 01 - 2020"
$ws.Range("E28").Value = "Done"
$ws.Range("F28").Value = 2

# 4. Extend the running-total formula (previously row 28, now row 29)
#    so it includes the new row.
$ws.Range("F29").Formula = "=SUM(F7:F28)"

# 5. Move the active selection the same way Excel would after inserting
#    a row above the previously selected cell (F29 -> F30).
$ws.Range("F30").Select()
